# Update the column-B header text from "Release Definition Name" to "Release Name".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Release Name"

# Move/save the active selection to D2 (matches the saved sheet view state).
$ws.Range("D2").Select()
